$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The language rows (originally A2:B23, sorted alphabetically) need to be
# reordered in descending order of their 1987 value, and the two rows whose
# value is 0 (Russian, Uzbek) need to be dropped entirely, shrinking the
# table from A1:B23 down to A1:B21.

# Rewrite rows 2..21 in the new (value-descending, zero rows excluded) order.
$ws.Range("A2").Value = "English"
$ws.Range("B2").Value = 30.07990767484222
$ws.Range("A3").Value = "Spanish"
$ws.Range("B3").Value = 8.918494117142496
$ws.Range("A4").Value = "Japanese"
$ws.Range("B4").Value = 8.685536830794122
$ws.Range("A5").Value = "German"
$ws.Range("B5").Value = 7.135445880409627
$ws.Range("A6").Value = "Arabic"
$ws.Range("B6").Value = 4.91890336888511
$ws.Range("A7").Value = "Chinese"
$ws.Range("B7").Value = 4.785074232759748
$ws.Range("A8").Value = "Portuguese"
$ws.Range("B8").Value = 4.640464075484049
$ws.Range("A9").Value = "French"
$ws.Range("B9").Value = 4.35655930865532
$ws.Range("A10").Value = "Italian"
$ws.Range("B10").Value = 4.284708550335647
$ws.Range("A11").Value = "Malay-Indonesian"
$ws.Range("B11").Value = 2.091945164264509
$ws.Range("A12").Value = "Dutch"
$ws.Range("B12").Value = 1.844735849534926
$ws.Range("A13").Value = "Persian"
$ws.Range("B13").Value = 1.522795385668429
$ws.Range("A14").Value = "Turkish"
$ws.Range("B14").Value = 1.446242412055626
$ws.Range("A15").Value = "Polish"
$ws.Range("B15").Value = 1.039762023382776
$ws.Range("A16").Value = "Korean"
$ws.Range("B16").Value = 0.9993138854554894
$ws.Range("A17").Value = "Urdu"
$ws.Range("B17").Value = 0.7514094606390961
$ws.Range("A18").Value = "Thai"
$ws.Range("B18").Value = 0.6898399711197661
$ws.Range("A19").Value = "Swedish"
$ws.Range("B19").Value = 0.6599320645741824
$ws.Range("A20").Value = "Bengali"
$ws.Range("B20").Value = 0.4076661647032911
$ws.Range("A21").Value = "Vietnamese"
$ws.Range("B21").Value = 0.2658182749621074

# Drop the two now-unused trailing rows (formerly Uzbek / Vietnamese) so the
# sheet shrinks from A1:B23 to A1:B21. Delete the lower row first so the
# upper row number stays valid.
$ws.Rows(23).Delete()
$ws.Rows(22).Delete()
